$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 11-51 shift down by one position because a new coin (OKB) enters the
# ranking at #9 (row 11); the former bottom row (Aave) drops off the list.
# Coin name + link (column B/C) are updated for every row that shifted,
# and every row also gets a refreshed Price (D) and Volume(1h) (E).
$coinRows = @(
    @{ Row = 11; Coin = "OKB"; Link = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; Price = "41.61"; Volume = "  +0.04%  " }
    @{ Row = 12; Coin = "Polkadot"; Link = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; Price = "6.235"; Volume = "  +0.67%  " }
    @{ Row = 14; Coin = "Solana"; Link = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; Price = "20.68"; Volume = "  +1.14%  " }
    @{ Row = 15; Coin = "Chainlink"; Link = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; Price = "7.267"; Volume = "  +1.21%  " }
    @{ Row = 16; Coin = "BinanceUSD"; Link = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; Price = "1.008"; Volume = "  +0.38%  " }
    @{ Row = 17; Coin = "ShibaInu"; Link = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; Price = "0.00001102"; Volume = "  +0.71%  " }
    @{ Row = 18; Coin = "Litecoin"; Link = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; Price = "91.21"; Volume = "  +0.83%  " }
    @{ Row = 19; Coin = "TRON"; Link = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; Price = "0.06680"; Volume = "  +1.05%  " }
    @{ Row = 20; Coin = "Avalanche"; Link = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; Price = "17.81"; Volume = "  +0.32%  " }
    @{ Row = 21; Coin = "Dai"; Link = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; Price = "1.006"; Volume = "  +0.36%  " }
    @{ Row = 22; Coin = "Uniswap"; Link = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; Price = "6.053"; Volume = "  +0.73%  " }
    @{ Row = 23; Coin = "WrappedBTC"; Link = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; Price = "28.238.98"; Volume = "  +0.76%  " }
    @{ Row = 24; Coin = "Cosmos"; Link = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; Price = "11.15"; Volume = "  +1.03%  " }
    @{ Row = 25; Coin = "Toncoin"; Link = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; Price = "2.264"; Volume = "  +1.65%  " }
    @{ Row = 26; Coin = "WrappedliquidstakedEther2.0"; Link = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; Price = "2.093.46"; Volume = "  +0.85%  " }
    @{ Row = 27; Coin = "Monero"; Link = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; Price = "160.09"; Volume = "  +1.45%  " }
    @{ Row = 28; Coin = "LidoDAOToken"; Link = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; Price = "2.485"; Volume = "  -2.79%  " }
    @{ Row = 29; Coin = "EthereumClassic"; Link = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; Price = "20.77"; Volume = "  +1.91%  " }
    @{ Row = 30; Coin = "BitcoinCash"; Link = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; Price = "125.32"; Volume = "  +0.58%  " }
    @{ Row = 31; Coin = "Stellar"; Link = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; Price = "0.1062"; Volume = "  +0.36%  " }
    @{ Row = 32; Coin = "ImmutableX"; Link = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; Price = "1.040"; Volume = "  +0.57%  " }
    @{ Row = 33; Coin = "Filecoin"; Link = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; Price = "5.849"; Volume = "  +4.51%  " }
    @{ Row = 34; Coin = "HuobiToken"; Link = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; Price = "3.616"; Volume = "  +0.21%  " }
    @{ Row = 35; Coin = "FraxShare"; Link = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; Price = "9.668"; Volume = "  +0.72%  " }
    @{ Row = 36; Coin = "VeChain"; Link = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; Price = "0.02452"; Volume = "  +1.72%  " }
    @{ Row = 37; Coin = "Hedera"; Link = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; Price = "0.06569"; Volume = "  +0.51%  " }
    @{ Row = 38; Coin = "Algorand"; Link = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; Price = "0.2190"; Volume = "  +0.90%  " }
    @{ Row = 39; Coin = "ARBITRUM"; Link = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; Price = "1.201"; Volume = "  -0.05%  " }
    @{ Row = 40; Coin = "TheSandbox"; Link = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; Price = "0.6501"; Volume = "  +1.42%  " }
    @{ Row = 41; Coin = "InternetComputer(DFINITY)"; Link = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; Price = "4.999"; Volume = "  +2.72%  " }
    @{ Row = 42; Coin = "TrustWalletToken"; Link = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; Price = "1.225"; Volume = "  -1.39%  " }
    @{ Row = 43; Coin = "Aptos"; Link = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; Price = "11.34"; Volume = "  +1.12%  " }
    @{ Row = 44; Coin = "Decentraland"; Link = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; Price = "0.6147"; Volume = "  +0.86%  " }
    @{ Row = 45; Coin = "EnergySwap"; Link = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; Price = "13.07"; Volume = "  +0.44%  " }
    @{ Row = 46; Coin = "WEMIXTOKEN"; Link = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; Price = "1.286"; Volume = "  +0.58%  " }
    @{ Row = 47; Coin = "PancakeSwap"; Link = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; Price = "3.680"; Volume = "  +0.61%  " }
    @{ Row = 48; Coin = "NEARProtocol"; Link = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; Price = "2.022"; Volume = "  +2.28%  " }
    @{ Row = 49; Coin = "EOS"; Link = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; Price = "1.231"; Volume = "  +2.30%  " }
    @{ Row = 50; Coin = "Quant"; Link = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; Price = "120.86"; Volume = "  +0.06%  " }
    @{ Row = 51; Coin = "Cronos"; Link = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; Price = "0.06916"; Volume = "  +1.09%  " }
)

foreach ($item in $coinRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $item.Coin
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $item.Link
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $item.Price
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $item.Volume
}

# Rows whose coin stayed in place still get refreshed Price/Volume figures.
$priceRows = @(
    @{ Row = 2; Price = "28.204.43"; Volume = "  +0.67%  " }
    @{ Row = 3; Price = "1.880.94"; Volume = "  +1.19%  " }
    @{ Row = 4; Price = "1.008"; Volume = "  +0.41%  " }
    @{ Row = 5; Price = "314.59"; Volume = "  +0.88%  " }
    @{ Row = 6; Price = "1.007"; Volume = "  +0.35%  " }
    @{ Row = 7; Price = "0.5142"; Volume = "  +1.02%  " }
    @{ Row = 8; Price = "0.3921"; Volume = "  +2.87%  " }
    @{ Row = 9; Price = "0.08371"; Volume = "  +1.52%  " }
    @{ Row = 10; Price = "1.122"; Volume = "  +1.14%  " }
    @{ Row = 13; Price = "1.889.38"; Volume = "  +1.49%  " }
)

foreach ($item in $priceRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $item.Price
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $item.Volume
}
